# Adds a new Q&A block to the end of the "Project note" document, right
# after the "//Why it work? ... <a>!!" paragraph, matching the commit
# "fixing the search .map".
#
# The new content (9 paragraphs, all sharing the same 100-twip first-line
# indent used throughout this section of the notes):
#   (blank)
#   (blank)
#   Q: Where should I put that code?
#   A: I need to put it in the innermost component, ...
#   But it in the Bookshelf component
#   (blank)
#   (blank)
#   Q: Why use promises?
#   A: To catch error
#
# The document's hidden "_GoBack" bookmark (Word's "last edit" marker) is
# always attached to the trailing edge of the document, so it has to move
# from the old last paragraph to the new one.

$d = $word.ActiveDocument

# The "_GoBack" bookmark is hidden (leading underscore) so it never shows
# up in $d.Bookmarks.Count, but it can still be fetched by name and removed
# so it doesn't linger on the paragraph that is no longer last.
$goBack = $d.Bookmarks.Item("_GoBack")
if ($goBack -ne $null) {
    $goBack.Delete()
}

# Insertion point: right after the existing text of the last paragraph,
# but before its end-of-paragraph mark.
$lastPara = $d.Paragraphs.Last
$insertAt = $lastPara.Range.End - 1
$target = $d.Range($insertAt, $insertAt)

# Insert the new paragraphs as a WordprocessingML fragment so each one gets
# exactly the formatting (and, for the final one, the relocated bookmark)
# that the target revision expects - empty paragraphs stay run-less, just
# like the existing blank paragraphs elsewhere in this document.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="100"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="100"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="100"/></w:pPr><w:r><w:t>Q: Where should I put that code?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="100"/></w:pPr><w:r><w:t>A: I need to put it in the innermost component, because in this case event propagates up. If I put it into App.js, it’s gonna cascade down to all components, including Search. Can’t put it in Search because you don’t want the button to show up on the Search page.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="100"/></w:pPr><w:r><w:t>But it in the Bookshelf component</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="100"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="100"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="100"/></w:pPr><w:r><w:t>Q: Why use promises?</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="100"/></w:pPr><w:r><w:t>A: To catch error</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$target.InsertXML($xml)
